$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 26772.94921875001
$ws.Range("D3").Value = 1777.34375
$ws.Range("E3").Value = 1173.046875
$ws.Range("F3").Value = 820.3125

$ws.Range("C4").Value = 45896.48437500001
$ws.Range("D4").Value = 3046.875
$ws.Range("E4").Value = 2010.9375
$ws.Range("F4").Value = 1406.25

$ws.Range("C5").Value = 57370.60546875001
$ws.Range("D5").Value = 3808.59375
$ws.Range("E5").Value = 2513.671875
$ws.Range("F5").Value = 1757.8125

$ws.Range("C6").Value = 61195.31250000001
$ws.Range("D6").Value = 4062.5
$ws.Range("E6").Value = 2681.25
$ws.Range("F6").Value = 1875

$ws.Range("C7").Value = 57370.60546875001
$ws.Range("D7").Value = 3808.59375
$ws.Range("E7").Value = 2513.671875
$ws.Range("F7").Value = 1757.8125

$ws.Range("C8").Value = 45896.48437500001
$ws.Range("D8").Value = 3046.875
$ws.Range("E8").Value = 2010.9375
$ws.Range("F8").Value = 1406.25

$ws.Range("C9").Value = 26772.94921875001
$ws.Range("D9").Value = 1777.34375
$ws.Range("E9").Value = 1173.046875
$ws.Range("F9").Value = 820.3125

$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
